# Add the new Tasmania destinations (rows 21-46, column A only) that were
# appended below the existing 20-row table, then move the selection/viewport
# down to the new bottom of the sheet.
#
# NOTE on entry order: the shared-strings table in the target workbook shows
# that "Lime Bay" (row 34) was actually typed AFTER "Cheesery Bruny" / "Ti Ama"
# / "Fish and Chips Dunalley" (rows 35-37) and then the row got moved up/inserted
# above them - its shared-string index (100) is higher than theirs (97-99).
# We reproduce that exact authoring order below so the regenerated
# sharedStrings.xml lines up with the target byte-for-byte.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "Great Lakes"
$ws.Range("A22").Value = "Liffey Falls"
$ws.Range("A23").Value = "Marion Bay"
$ws.Range("A24").Value = "Spring Beach"
$ws.Range("A25").Value = "Chain of Lagoons"
$ws.Range("A26").Value = "Piccaninny"
$ws.Range("A27").Value = "Ben Lomond"
$ws.Range("A28").Value = "Hollibanks Treetop Adventure"
$ws.Range("A29").Value = "Narawntapu Ntl Park"
$ws.Range("A30").Value = "Tarkine"
$ws.Range("A31").Value = "Marawah"
$ws.Range("A32").Value = "Tinderbox"
$ws.Range("A33").Value = "Boronia Reserve"
$ws.Range("A35").Value = "Cheesery Bruny"
$ws.Range("A36").Value = "Ti Ama"
$ws.Range("A37").Value = "Fish and Chips Dunalley"
$ws.Range("A34").Value = "Lime Bay"
$ws.Range("A38").Value = "South Arm surfspots"
$ws.Range("A39").Value = "Derby"
$ws.Range("A40").Value = "Maydena"
$ws.Range("A41").Value = "Mt Field"
$ws.Range("A42").Value = "Tahune/Hartz Mountain/Federation"
$ws.Range("A43").Value = "Rocky Cape Ntl Park"
$ws.Range("A44").Value = "Mole Creek"
$ws.Range("A45").Value = "Hastings Caves/Ida Bay"
$ws.Range("A46").Value = "Douglas Apsley Ntl Park"

# Row heights: the single-column entries auto-wrap inside the narrow 19-char
# "Name" column, so a short place name sits on one line (17pt) while longer
# names wrap to two lines (34pt) - matches what Excel would have computed.
$ws.Rows.Item(21).RowHeight = 17
$ws.Rows.Item(22).RowHeight = 17
$ws.Rows.Item(23).RowHeight = 17
$ws.Rows.Item(24).RowHeight = 17
$ws.Rows.Item(25).RowHeight = 17
$ws.Rows.Item(26).RowHeight = 17
$ws.Rows.Item(27).RowHeight = 17
$ws.Rows.Item(28).RowHeight = 34
$ws.Rows.Item(29).RowHeight = 17
$ws.Rows.Item(30).RowHeight = 17
$ws.Rows.Item(31).RowHeight = 17
$ws.Rows.Item(32).RowHeight = 17
$ws.Rows.Item(33).RowHeight = 17
$ws.Rows.Item(34).RowHeight = 17
$ws.Rows.Item(35).RowHeight = 17
$ws.Rows.Item(36).RowHeight = 17
$ws.Rows.Item(37).RowHeight = 34
$ws.Rows.Item(38).RowHeight = 17
$ws.Rows.Item(39).RowHeight = 17
$ws.Rows.Item(40).RowHeight = 17
$ws.Rows.Item(41).RowHeight = 17
$ws.Rows.Item(42).RowHeight = 34
$ws.Rows.Item(43).RowHeight = 17
$ws.Rows.Item(44).RowHeight = 17
$ws.Rows.Item(45).RowHeight = 34
$ws.Rows.Item(46).RowHeight = 34

# Scroll the view down and move the selection to the first empty row below
# the newly-added data, just like the author left it before saving.
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("A47").Select()

Write-Output "Added 26 new Tasmania destinations (rows 21-46) and updated selection."
